# TEST 1 - ETS 2024 - passages: remove the "listening" section rows and
# rename the remaining "reading" passage codes to use a hyphen after
# "passage" (e.g. "passage131-134" -> "passage-131-134").

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-6 hold the "listening" section entries; deleting them shifts the
# "reading" rows (formerly 7-25) up to rows 2-20, and shrinks the trailing
# block of empty filler rows by the same 5 rows (1000 -> 995).
$ws.Rows("2:6").Delete()

# Rename passage codes in column B (now rows 2-20) to include a hyphen
# right after "passage".
for ($r = 2; $r -le 20; $r++) {
    $cell = $ws.Cells.Item($r, 2)
    $val = $cell.Value2
    if ($val -like "passage*") {
        $cell.Value = "passage-" + $val.Substring(7)
    }
}
